# -----------------------------------------------------------------------
# Adds the new quarter "2022-Q3" fund-holdings sheet and its row in the
# "总计" (totals) summary sheet, pushing everything older down by one.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" sheet: insert a fresh row 2 for "2022-Q3"; the existing
#    quarter rows shift down to row 3, 4, 5, ... (their A-column running
#    index already equals row-2, so it stays correct after the shift).
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)
$total.Rows.Item(2).Insert()

# Insert() doesn't copy the index-column style (s="2") onto the new A2;
# grab it from A3 (just shifted down, still correctly styled).
$total.Cells.Item(3,1).Copy()
$total.Cells.Item(2,1).PasteSpecial(-4122)   # xlPasteFormats

# Insert() also leaves a stray style on the new B2:D2; clear it so those
# cells are plain, like the rest of the data rows.
$total.Range("B2:D2").ClearFormats()

$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q3"
$total.Cells.Item(2,3).Value = 15
$total.Cells.Item(2,4).Value = 5.58

# ---------------------------------------------------------------------
# 2) New "2022-Q3" holdings sheet, inserted right after "总计" (i.e.
#    before the sheet currently named "2022-Q2"). Clone "2022-Q2" so
#    the new sheet starts with identical column widths/header styling,
#    then overwrite its title and all data cells.
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item(2)
$q2.Copy($q2)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# The cloned sheet carries 17 data rows (2022-Q2 had 17 funds); 2022-Q3
# only has 15, so drop the two extra trailing rows first.
$q3.Rows.Item(18).Delete()
$q3.Rows.Item(17).Delete()

# Columns: A idx(n) | B code(text) | C name(text) | D size(text)
#        | E position(text) | F pct(text) | G value(text) | H rank(n)
$q3Data = @(
    @(0, "540008", "汇丰晋信低碳先锋股票A", "79.15", "94.63", "4.36", "3.4509", 10),
    @(1, "011578", "汇丰晋信核心成长混合A", "21.81", "94.20", "3.76", "0.8201", 9),
    @(2, "001643", "汇丰晋信智造先锋股票A", "15.68", "94.47", "3.72", "0.5833", 9),
    @(3, "001644", "汇丰晋信智造先锋股票C", "8.27", "94.47", "3.72", "0.3076", 9),
    @(4, "011579", "汇丰晋信核心成长混合C", "4.19", "94.20", "3.76", "0.1575", 9),
    @(5, "013511", "汇丰晋信低碳先锋股票C", "2.72", "94.63", "4.36", "0.1186", 10),
    @(6, "000963", "兴业多策略灵活配置混合", "1.63", "83.65", "4.13", "0.0673", 3),
    @(7, "008602", "方正富邦新兴成长混合A", "1.23", "86.03", "3.90", "0.0480", 7),
    @(8, "008082", "国寿安保研究精选混合A", "0.37", "92.74", "2.92", "0.0108", 10),
    @(9, "014305", "华泰柏瑞中证500指数增强A", "2.20", "34.78", "0.41", "0.0090", 8),
    @(10, "008083", "国寿安保研究精选混合C", "0.16", "92.74", "2.92", "0.0047", 10),
    @(11, "012415", "德邦上证G60综指增强A", "0.09", "92.83", "2.29", "0.0021", 8),
    @(12, "008603", "方正富邦新兴成长混合C", "0.03", "86.03", "3.90", "0.0012", 7),
    @(13, "014306", "华泰柏瑞中证500指数增强C", "0.07", "34.78", "0.41", "0.0003", 8),
    @(14, "012416", "德邦上证G60综指增强C", "0.01", "92.83", "2.29", "0.0002", 8)
)

$r = 2
foreach ($row in $q3Data) {
    $q3.Cells.Item($r,1).Value = $row[0]
    $q3.Cells.Item($r,2).Value = "'" + $row[1]
    $q3.Cells.Item($r,3).Value = $row[2]
    $q3.Cells.Item($r,4).Value = "'" + $row[3]
    $q3.Cells.Item($r,5).Value = "'" + $row[4]
    $q3.Cells.Item($r,6).Value = "'" + $row[5]
    $q3.Cells.Item($r,7).Value = "'" + $row[6]
    $q3.Cells.Item($r,8).Value = $row[7]
    $r = $r + 1
}

# Restore the selection/active sheet to "总计" (matches the original
# workbook's activeTab, left untouched by this edit otherwise).
$wb.Worksheets.Item(1).Activate()
